$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Recorded By" (column G) lists users as a comma-separated string.
# Reorder each list so that email-like entries come first, followed by
# the remaining (non-email) entries such as "System"/"system", each
# group keeping its original relative order.
for ($row = 2; $row -le 157; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $val = $cell.Text

    if ([string]::IsNullOrEmpty($val)) { continue }

    $parts = $val -split ',\s*'
    $emails = @()
    $others = @()
    foreach ($p in $parts) {
        if ($p -like "*@*") {
            $emails += $p
        } else {
            $others += $p
        }
    }
    $newParts = $emails + $others
    $newVal = [string]::Join(", ", $newParts)

    if ($newVal -ne $val) {
        $cell.Value = $newVal
    }
}
